# --- inValidLoginData: drop the old selection, select A1:B1 instead -------
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
[void]$ws1.Range("A1:B1").Select()

# --- add the new "validLoginData" sheet right after it --------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "validLoginData"

# Header row.
$ws2.Range("A1").Value = "userName"
$ws2.Range("B1").Value = "password"

# Valid credentials row.
$ws2.Range("A2").Value = "Admin"
$ws2.Range("B2").Value = "admin123"

# Column widths (characters), matching the source workbook's layout.
$ws2.Columns.Item(1).ColumnWidth = 12.42
$ws2.Columns.Item(2).ColumnWidth = 12.25

# Header formatting: reuse the same bold/centered/wrapped look already
# used for the header row on inValidLoginData, via a plain format copy
# (keeps the shared style table tidy instead of minting new xfs).
[void]$ws1.Range("A1:B1").Copy()
[void]$ws2.Range("A1:B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws2.Rows.Item(1).RowHeight = 30

# Leave the cursor on B2 -- the active cell on the now-active new sheet.
[void]$ws2.Range("B2").Select()
